$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure each written cell keeps a Text format, matching the source
# workbook (all cells in this sheet are stored as inline strings),
# so numeric-looking values (e.g. "49.10", "16.50") are not silently
# normalized to "49.1" / "16.5" by automatic type inference.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.366.07'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.568.07'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.88'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3749'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.10'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07557'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.94'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.948'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.911'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.566.17'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.92'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06745'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.50'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.175'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.369.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.378'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.700'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.82'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.036'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.46'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.741.95'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.033'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9849'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.984'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.431'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +11.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08456'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02485'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2282'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06463'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.387'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.13'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.91'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.803'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5900'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.061'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '124.63'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07324'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.80%  '
